$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5 (ALC)
$ws.Range("H5").Value = 247.33333
$ws.Range("I5").Value = 247.33333
$ws.Range("K5").Value = 247.33333
$ws.Range("M5").Value = -132.33333

# Row 18 (ALC)
$ws.Range("H18").Value = 9640.263000000001
$ws.Range("I18").Value = 8715.588
$ws.Range("K18").Value = 8715.588
$ws.Range("M18").Value = -8431.588

# Row 19 (ALC)
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()

# Row 98 (ALC)
$ws.Range("H98").Value = 2004.3334
$ws.Range("I98").Value = 1725.8334
$ws.Range("K98").Value = 1725.8334
$ws.Range("M98").Value = -227.8334

# Row 112 (ALC)
$ws.Range("H112").Value = 7292.125
$ws.Range("J112").Value = 10239.143
$ws.Range("L112").Value = 30717.429
$ws.Range("N112").Value = -32933.429

# Row 122 (ALC)
$ws.Range("H122").Value = 2004.3334
$ws.Range("I122").Value = 1725.8334
$ws.Range("K122").Value = 5177.5002
$ws.Range("M122").Value = -2727.5002

# Row 137 (ALC)
$ws.Range("H137").Value = 2178
$ws.Range("J137").Value = 2998.5
$ws.Range("L137").Value = 8995.5
$ws.Range("N137").Value = -14095.5

# Row 138 (ALC)
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (ARM)
$ws.Range("H45").Value = 6099.8
$ws.Range("I45").Value = 6099.8
$ws.Range("K45").Value = 6099.8
$ws.Range("M45").Value = -5722.8

# Row 102 (ARM)
$ws.Range("H102").Value = 1497.25
$ws.Range("I102").Value = 1497.25
$ws.Range("K102").Value = 1497.25
$ws.Range("M102").Value = 124.75

$ws = $wb.Worksheets.Item("CRP")
# Row 8 (CRP)
$ws.Range("H8").Value = 2300
$ws.Range("I8").Value = 2500
$ws.Range("J8").Value = 2200
$ws.Range("K8").Value = 2500
$ws.Range("L8").Value = 2200
$ws.Range("M8").Value = -2360
$ws.Range("N8").Value = -2480

# Row 22 (CRP)
$ws.Range("H22").Value = 679.1667
$ws.Range("I22").Value = 718.2727
$ws.Range("J22").Value = 249
$ws.Range("K22").Value = 718.2727
$ws.Range("L22").Value = 249
$ws.Range("M22").Value = -368.2727
$ws.Range("N22").Value = -949

# Row 62 (CRP)
$ws.Range("H62").Value = 8400.799999999999
$ws.Range("I62").Value = 10668
$ws.Range("K62").Value = 10668
$ws.Range("M62").Value = -10044

# Row 65 (CRP)
$ws.Range("H65").Value = 8400.799999999999
$ws.Range("I65").Value = 10668
$ws.Range("K65").Value = 53340
$ws.Range("M65").Value = -50220

$ws = $wb.Worksheets.Item("CUL")
# Row 33 (CUL)
$ws.Range("H33").Value = 9
$ws.Range("I33").Value = 9
$ws.Range("J33").Value = 9
$ws.Range("K33").Value = 54
$ws.Range("L33").Value = 54
$ws.Range("M33").Value = 229
$ws.Range("N33").Value = -620

# Row 62 (CUL)
$ws.Range("H62").Value = 18000
$ws.Range("I62").Value = 18000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 54000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -53314
$ws.Range("N62").ClearContents()

# Row 65 (CUL)
$ws.Range("H65").Value = 18000
$ws.Range("I65").Value = 18000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 162000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -158568
$ws.Range("N65").ClearContents()

# Row 86 (CUL)
$ws.Range("H86").Value = 15309.167
$ws.Range("I86").Value = 7950
$ws.Range("J86").Value = 22668.334
$ws.Range("K86").Value = 23850
$ws.Range("L86").Value = 68005.00199999999
$ws.Range("M86").Value = -22664
$ws.Range("N86").Value = -70377.00199999999

# Row 89 (CUL)
$ws.Range("H89").Value = 15309.167
$ws.Range("I89").Value = 7950
$ws.Range("J89").Value = 22668.334
$ws.Range("K89").Value = 71550
$ws.Range("L89").Value = 204015.006
$ws.Range("M89").Value = -65622
$ws.Range("N89").Value = -215871.006

# Row 122 (CUL)
$ws.Range("H122").Value = 1499.3334
$ws.Range("J122").Value = 1499.3334
$ws.Range("L122").Value = 13494.0006
$ws.Range("N122").Value = -18394.0006

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (GSM)
$ws.Range("H80").Value = 200006
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 200006
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 200006
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -202002

# Row 83 (GSM)
$ws.Range("H83").Value = 200006
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 200006
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 1000030
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -1010014

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (LTW)
$ws.Range("H22").Value = 3000.3333
$ws.Range("I22").Value = 3000.5
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 3000.5
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -2705.5
$ws.Range("N22").Value = -3590

# Row 27 (LTW)
$ws.Range("H27").Value = 3000.3333
$ws.Range("I27").Value = 3000.5
$ws.Range("J27").Value = 3000
$ws.Range("K27").Value = 3000.5
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -2893.5
$ws.Range("N27").Value = -3214

# Row 46 (LTW)
$ws.Range("H46").Value = 3049.9
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 3437.375
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 3437.375
$ws.Range("M46").Value = -1312
$ws.Range("N46").Value = -3813.375

# Row 82 (LTW)
$ws.Range("H82").Value = 2493.625
$ws.Range("I82").Value = 2135.5715
$ws.Range("K82").Value = 2135.5715
$ws.Range("M82").Value = -1774.5715

# Row 85 (LTW)
$ws.Range("H85").Value = 2493.625
$ws.Range("I85").Value = 2135.5715
$ws.Range("K85").Value = 2135.5715
$ws.Range("M85").Value = -887.5715

# Row 93 (LTW)
$ws.Range("H93").Value = 20999.5
$ws.Range("I93").Value = 20999.5
$ws.Range("K93").Value = 20999.5
$ws.Range("M93").Value = -19751.5

# Row 94 (LTW)
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# Row 100 (LTW)
$ws.Range("H100").Value = 3875
$ws.Range("I100").Value = 4166.6665
$ws.Range("K100").Value = 4166.6665
$ws.Range("M100").Value = -3625.6665
